$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 173. This shifts the former rows
# 173-182 down to 175-184 (their contents remain unchanged).
$ws.Range("A173:A174").EntireRow.Insert()

# --- New row 173 ---
$ws.Range("A173").Value2 = 9
$ws.Range("B173").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C173").Value2 = "Metropolitana"
$ws.Range("D173").Value2 = 45267
$ws.Range("D173").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E173").Value2 = 13
$ws.Range("F173").Value2 = "Fruta"
$ws.Range("G173").Value2 = 100101
$ws.Range("H173").Value2 = "Berries"
$ws.Range("I173").Value2 = 100101004
$ws.Range("J173").Value2 = "Frambuesa"
$ws.Range("K173").Value2 = "Sin especificar"
$ws.Range("L173").Value2 = "Especial"
$ws.Range("M173").Value2 = 500
$ws.Range("N173").Value2 = 8000
$ws.Range("O173").Value2 = 8000
$ws.Range("P173").Value2 = 8000
$ws.Range("Q173").Value2 = "$/bandeja 2 kilos"
$ws.Range("R173").Value2 = "Provincia de Curicó"
$ws.Range("S173").Value2 = 4000
$ws.Range("T173").Value2 = 2

# --- New row 174 ---
$ws.Range("A174").Value2 = 9
$ws.Range("B174").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C174").Value2 = "Metropolitana"
$ws.Range("D174").Value2 = 45267
$ws.Range("D174").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E174").Value2 = 13
$ws.Range("F174").Value2 = "Fruta"
$ws.Range("G174").Value2 = 100101
$ws.Range("H174").Value2 = "Berries"
$ws.Range("I174").Value2 = 100101004
$ws.Range("J174").Value2 = "Frambuesa"
$ws.Range("K174").Value2 = "Sin especificar"
$ws.Range("L174").Value2 = "Primera"
$ws.Range("M174").Value2 = 250
$ws.Range("N174").Value2 = 7000
$ws.Range("O174").Value2 = 7000
$ws.Range("P174").Value2 = 7000
$ws.Range("Q174").Value2 = "$/bandeja 2 kilos"
$ws.Range("R174").Value2 = "Provincia de Curicó"
$ws.Range("S174").Value2 = 3500
$ws.Range("T174").Value2 = 2
